$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Capture the existing hyperlinks (cell address + target) in their current collection order ---
$links = @()
foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    $links += ,@($addr, $hl.Address)
}

# --- Drop all hyperlink objects; they will be rebuilt (with corrected addresses) after the row shift ---
$ws.Hyperlinks.Delete()

# --- Delete the entire "ATIPPAL UAT" row (row 12); https://atippal-uat.ssc-spc.gc.ca / ATIPPAL UAT ---
$ws.Rows.Item(12).Delete()

# --- Recompute each captured hyperlink's new cell address: drop row 12, shift rows below it up by one ---
function Shift-RowInAddress($addr) {
    $parts = $addr -split '\$'
    $col = $parts[1]
    $row = [int]$parts[2]
    if ($row -eq 12) {
        return $null
    } elseif ($row -gt 12) {
        $row = $row - 1
    }
    return "$col" + "$row"
}

$newCells = @()
foreach ($l in $links) {
    $oldAddr = $l[0]
    $target = $l[1]
    $newCell = Shift-RowInAddress $oldAddr
    if ($newCell -ne $null) {
        $rng = $ws.Range($newCell)
        $ws.Hyperlinks.Add($rng, $target) | Out-Null
        $newCells += $newCell
    }
}

# --- Re-adding hyperlinks resets formatting; restore the original "Hyperlink" cell style ---
foreach ($c in $newCells) {
    $ws.Range($c).Style = "Hyperlink"
}

# --- Match the selection left on the sheet after the row deletion ---
$ws.Rows.Item(12).Select()
